$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 14, shifting existing rows 14-16 down to 15-17.
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with a duplicate of the old row 13
# data (dated 2022-09-09), then bump row 13's date forward to the new
# weekly reading (2022-09-20).
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44813
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100112036
$ws.Range("G14").Value = "Caigua"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 20
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 20000
$ws.Range("N14").Value = "$/caja 15 kilos"
$ws.Range("O14").Value = "Región de Arica y Parinacota"
$ws.Range("P14").Value = 1333
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = "Hortaliza"

# Row 13 keeps all its original values except the date, which advances to
# the new week.
$ws.Range("D13").Value = 44824
